# The workbook was previously saved by a non-Excel tool (Apache POI) with
# placeholder view/selection state and a stray blank "ReferenceNumber"
# column (E) that never received data. Re-touch it the way Excel itself
# would when a user opens the sheet, clears the unused column, and leaves
# the selection sitting on the real data block before saving again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E (ReferenceNumber) only ever held its header; every data row's
# cell is blank. Clear those out so the cells disappear from the saved
# sheet instead of lingering as empty <c> elements.
$ws.Range("E2:E151").ClearContents() | Out-Null

# Leave the sheet's selection on the data body (matches the file as
# reopened/resaved), with the sheet itself the active tab of the window.
$ws.Activate() | Out-Null
$ws.Range("A5:G151").Select() | Out-Null
